# working_hours.xlsx - add a missed time entry (2014-05-06, 17:30-18:30)
#
# A new data row is inserted right before the current row 53 (the blank
# spacer row that precedes the "sum [min]" row). This pushes the spacer
# row and the three summary rows (sum [min] / sum [h] / sum [working
# weeks]) down by one row, and Excel's formula-reference adjustment
# naturally extends:
#   - the shared formulas in columns F/G (F28:F52 -> F28:F53, etc.)
#   - the SUM(F2:F53) total -> SUM(F2:F54)
# An existing entry's end time (E52) is also corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 53; everything from the old row 53
# downward (spacer row + the three summary rows) shifts down by one.
$ws.Rows.Item(53).Insert()

# Correct the end time of the entry on row 52 (was 0.75 / 18:00, now
# 0.6875 / 16:30). F52/G52 recompute automatically via the shared formula.
$ws.Range("E52").Value = 0.6875

# Populate the newly inserted row 53 with the missed entry.
$ws.Range("A53").Value = 2014
$ws.Range("B53").Value = 5
$ws.Range("C53").Value = 6
$ws.Range("D53").Value = 0.72916666666666663
$ws.Range("E53").Value = 0.77083333333333337
$ws.Range("F53").Formula = "=(E53-D53)*24*60"
$ws.Range("G53").Formula = "=F53/60"

# Match the author's final selection.
$ws.Range("E54").Select() | Out-Null
